# ST-63 Mockup Delete Ticket edits
$d = $word.ActiveDocument

# --- 1) Image run right after "So do luong du lieu:" gets <w:noProof/> ---
$d.Paragraphs(2).Range.NoProofing = 1

# --- 2) D2 paragraph: drop trailing period, append extra clause ---
$p = $d.Paragraphs(5)
$full = $d.Range($p.Range.Start, $p.Range.End)
$full.Text = "D2: Lấy danh sách các chuyến bay đã chọn ở D1 để xóa vé tự động, Các chuyến bay đó phải chưa được diễn ra"

# --- 3) Insert a brand-new "Buoc 5" step paragraph right after "Buoc 4" ---
$p4 = $d.Paragraphs(12)
$newPara = $p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(13)
$p5start = $d.Range($p5.Range.Start, $p5.Range.Start)
$p5start.InsertAfter("Bước 5: Kiểm tra qui định “Chuyến bay”")

# --- 4) Renumber the remaining steps (old 5..9 -> new 6..10) ---
$p = $d.Paragraphs(14)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = "Bước 6: Thay đổi tình trạng xóa ẩn"

$p = $d.Paragraphs(15)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = "Bước 7: Lưu xuống cở sở dữ liệu sau khi thay đổi tình trạng"

$p = $d.Paragraphs(16)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = "Bước 8: Hiển thị lại danh sách vé sau khi xóa"

$p = $d.Paragraphs(17)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = "Bước 9: Đóng kết nối cơ sở dữ liệu"

$p = $d.Paragraphs(18)
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Text = "Bước 10: Kết thúc"

# --- 5) Append a new "Luu y" paragraph after the "... xoa theo." explanation ---
$p = $d.Paragraphs(22)
$newPara2 = $p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(23)
$p2start = $d.Range($p2.Range.Start, $p2.Range.Start)
$p2start.InsertAfter("Lưu ý: Chỉ được xóa vé của các chuyến bay chưa diễn ra và không được xóa các chuyến bay mà không thể xóa như: chuyến bay đó đã có số lượng hành khách quá nửa của chuyến bay đó, ….")
